$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header update
$ws.Range("E1").Value = "strength (RMS)"

# Data rows (experiment, reactionTime, peakTime, difference, strength)
$data = @(
    @("April16_offset_PPI_18", 24.67, 30.67, 6, 103),
    @("April16_ASR_control", 19, 27.5, 8.5, 69.25),
    @("April16_offset_PPI_12", 24.8, 30, 5.2, 82.40000000000001),
    @("April16_gap_depth", 17.2, 30.4, 13.2, 108.4),
    @("April16_gap_duration_20", 18, 30, 12, 96.75),
    @("April16_gap_duration_8", 21.2, 26, 4.8, 97.40000000000001),
    @("April16_offset_PPI_20", 24, 30.8, 6.8, 52),
    @("April16_gap_duration_10", 16.8, 28.4, 11.6, 75.8),
    @("April16_gap_duration_50", 19.6, 30.4, 10.8, 90),
    @("April16_offset_PPI_50", 20.8, 30, 9.199999999999999, 79.2),
    @("April16_offset_PPI_14", 20.8, 30.8, 10, 86.2),
    @("April16_offset_PPI_16", 17, 30, 13, 104.25),
    @("April16_offset_PPI_6", 19, 30, 11, 93.25),
    @("April16_gap_duration_4", 21.5, 33, 11.5, 93.5),
    @("April16_tone_in_noise", 14.67, 30, 15.33, 107.33),
    @("April16_offset_PPI_10", 30, 32, 2, 85),
    @("April16_offset_PPI_4", 23, 30.5, 7.5, 75.5),
    @("April16_offset_PPI_8", 18.5, 30, 11.5, 64.25)
)

$row = 2
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
    $ws.Cells.Item($row, 5).Value = $r[4]
    $row++
}
